$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D9")
$r.VerticalAlignment = -4160
Write-Host "done"
